# Update the "dSF" (column F) values for several rows to reflect
# repulled data / recalculated means.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -7
$ws.Range("F3").Value = -6
$ws.Range("F6").Value = -10
$ws.Range("F7").Value = -5
$ws.Range("F8").Value = 2
$ws.Range("F11").Value = -8
$ws.Range("F12").Value = -2
